$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.378.85'
$ws.Range('E2').Value = '  +2.57%  '
$ws.Range('D3').Value = '2.349.79'
$ws.Range('E3').Value = '  +5.97%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.96'
$ws.Range('E5').Value = '  +5.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.48'
$ws.Range('E6').Value = '  +1.80%  '
$ws.Range('E7').Value = '  +2.99%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.635'
$ws.Range('E9').Value = '  +6.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.01'
$ws.Range('E10').Value = '  -1.19%  '
$ws.Range('E11').Value = '  +3.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.81'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('E13').Value = '  +6.87%  '
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.38'
$ws.Range('E15').Value = '  +9.32%  '
$ws.Range('D16').Value = '2.705.74'
$ws.Range('E16').Value = '  +6.17%  '
$ws.Range('D17').Value = '2.436.99'
$ws.Range('E17').Value = '  +8.83%  '
$ws.Range('D18').Value = '43.292.55'
$ws.Range('E18').Value = '  +2.39%  '
$ws.Range('E19').Value = '  +3.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.24'
$ws.Range('E20').Value = '  -2.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '75.52'
$ws.Range('E21').Value = '  +4.35%  '
$ws.Range('B22').Value = 'PancakeSwap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.45'
$ws.Range('E22').Value = '  -0.43%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.58'
$ws.Range('E23').Value = '  +12.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '254.86'
$ws.Range('E24').Value = '  +11.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.13'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.06'
$ws.Range('E26').Value = '  +4.18%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.20'
$ws.Range('E28').Value = '  +2.61%  '
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.37'
$ws.Range('E30').Value = '  +6.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.68'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0929'
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.07'
$ws.Range('E34').Value = '  +8.41%  '
$ws.Range('E35').Value = '  +5.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.99'
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.16'
$ws.Range('E37').Value = '  -4.05%  '
$ws.Range('E38').Value = '  +2.60%  '
$ws.Range('E39').Value = '  +2.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.70'
$ws.Range('E40').Value = '  +11.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.59'
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('E42').Value = '  +14.55%  '
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.82'
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  +4.08%  '
$ws.Range('E47').Value = '  +10.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '111.24'
$ws.Range('E48').Value = '  +7.84%  '
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.100'
$ws.Range('E50').Value = '  +3.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.462'
$ws.Range('E51').Value = '  +5.78%  '
